# Monthly update on files
# Updates the date column (A2:A26) from 2025-11-12 to 2025-11-30 and refreshes
# the metric columns (D2:O26) with the latest monthly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New report date for every data row
$newDate = Get-Date -Year 2025 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0

# D..O values per row (r=2..26), in column order D,E,F,G,H,I,J,K,L,M,N,O
$data = @(
    @(128,33,3,33,3,0,0,0,0,64,67,165),
    @(131,51,8,51,8,0,0,0,0,72,79,203),
    @(44,60,9,60,9,0,0,0,0,103,67,231),
    @(8,124,19,124,19,0,0,0,0,108,8,240),
    @(3,35,2,35,2,0,0,0,0,61,10,107),
    @(172,29,2,29,2,0,0,0,0,46,49,125),
    @(186,31,1,31,1,0,0,0,0,50,58,140),
    @(74,37,4,37,4,0,0,0,0,76,37,150),
    @(7,24,3,24,3,0,0,0,0,45,9,79),
    @(148,35,5,35,5,0,0,0,0,56,65,155),
    @(153,44,6,44,6,0,0,0,0,71,118,232),
    @(69,61,14,61,14,0,0,0,0,127,96,282),
    @(21,47,4,47,4,0,0,0,0,122,35,203),
    @(135,22,2,22,2,0,0,0,0,50,56,128),
    @(123,43,13,43,13,0,0,0,0,71,101,215),
    @(43,41,4,41,4,0,0,0,0,109,61,211),
    @(14,48,11,48,11,0,0,0,0,88,20,156),
    @(115,28,2,0,0,28,2,0,0,75,82,184),
    @(135,40,3,0,0,40,3,0,0,75,123,236),
    @(49,43,6,0,0,43,6,0,0,68,36,147),
    @(7,30,4,0,0,30,4,0,0,123,18,172),
    @(150,28,2,28,2,0,0,0,0,62,51,154),
    @(145,41,3,41,3,0,0,0,0,74,94,233),
    @(60,52,8,52,8,0,0,0,0,102,66,247),
    @(13,38,5,38,5,0,0,0,0,95,36,173)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newDate

    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        # Column D is the 4th column
        $col = 4 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
